# Auto-generated edit script: refresh market-price derived columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match latest Universalis pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 8: H8: 3350 -> 2514.2856; I8: 25.5 -> 20.2; J8: 9999 -> 8749.5; K8: 76.5 -> 60.59999999999999; L8: 29997 -> 26248.5; M8: 62.5 -> 78.40000000000001; N8: -30275 -> -26526.5
$ws.Range("H8").Value = 2514.2856
$ws.Range("I8").Value = 20.2
$ws.Range("J8").Value = 8749.5
$ws.Range("K8").Value = 60.59999999999999
$ws.Range("L8").Value = 26248.5
$ws.Range("M8").Value = 78.40000000000001
$ws.Range("N8").Value = -26526.5

# ALC row 39: H39: 29.666666 -> 26.714285; I39: 38.25 -> 32.2; J39: 12.5 -> 13; K39: 114.75 -> 96.60000000000001; L39: 37.5 -> 39; M39: 181.25 -> 199.4; N39: -629.5 -> -631
$ws.Range("H39").Value = 26.714285
$ws.Range("I39").Value = 32.2
$ws.Range("J39").Value = 13
$ws.Range("K39").Value = 96.60000000000001
$ws.Range("L39").Value = 39
$ws.Range("M39").Value = 199.4
$ws.Range("N39").Value = -631

# ALC row 50: H50: 0 -> 1000; J50: 0 -> 1000; L50: 0 -> 3000; N50: None -> -3950 | add: ['N50']
$ws.Range("H50").Value = 1000
$ws.Range("J50").Value = 1000
$ws.Range("L50").Value = 3000
$ws.Range("N50").Value = -3950

# ALC row 64: H64: 5000 -> 1666.3334; I64: 4000 -> 1999; J64: 6000 -> 1500; K64: 4000 -> 1999; L64: 6000 -> 1500; M64: -3752 -> -1751; N64: -6496 -> -1996
$ws.Range("H64").Value = 1666.3334
$ws.Range("I64").Value = 1999
$ws.Range("J64").Value = 1500
$ws.Range("K64").Value = 1999
$ws.Range("L64").Value = 1500
$ws.Range("M64").Value = -1751
$ws.Range("N64").Value = -1996

# ALC row 67: H67: 5000 -> 1666.3334; I67: 4000 -> 1999; J67: 6000 -> 1500; K67: 4000 -> 1999; L67: 6000 -> 1500; M67: -3142 -> -1141; N67: -7716 -> -3216
$ws.Range("H67").Value = 1666.3334
$ws.Range("I67").Value = 1999
$ws.Range("J67").Value = 1500
$ws.Range("K67").Value = 1999
$ws.Range("L67").Value = 1500
$ws.Range("M67").Value = -1141
$ws.Range("N67").Value = -3216

# ALC row 105: H105: 42096.145 -> 40640.363; J105: 42096.145 -> 40640.363; L105: 42096.145 -> 40640.363; N105: -49084.145 -> -47628.363
$ws.Range("H105").Value = 40640.363
$ws.Range("J105").Value = 40640.363
$ws.Range("L105").Value = 40640.363
$ws.Range("N105").Value = -47628.363

$ws = $wb.Worksheets.Item("ARM")
# ARM row 28: H28: 6000 -> 5999; I28: 6000 -> 5999; K28: 6000 -> 5999; M28: -5808 -> -5807
$ws.Range("H28").Value = 5999
$ws.Range("I28").Value = 5999
$ws.Range("K28").Value = 5999
$ws.Range("M28").Value = -5807

# ARM row 41: H41: 2284.1428 -> 2378.2; I41: 1758.4 -> 2073; J41: 3598.5 -> 3599; K41: 1758.4 -> 2073; L41: 3598.5 -> 3599; M41: -1344.4 -> -1659; N41: -4426.5 -> -4427
$ws.Range("H41").Value = 2378.2
$ws.Range("I41").Value = 2073
$ws.Range("J41").Value = 3599
$ws.Range("K41").Value = 2073
$ws.Range("L41").Value = 3599
$ws.Range("M41").Value = -1659
$ws.Range("N41").Value = -4427

# ARM row 88: H88: 2335 -> 1003.6667; I88: 2335 -> 1003.6667; K88: 2335 -> 1003.6667; M88: -1929 -> -597.6667
$ws.Range("H88").Value = 1003.6667
$ws.Range("I88").Value = 1003.6667
$ws.Range("K88").Value = 1003.6667
$ws.Range("M88").Value = -597.6667

# ARM row 91: H91: 2335 -> 1003.6667; I91: 2335 -> 1003.6667; K91: 2335 -> 1003.6667; M91: -931 -> 400.3333
$ws.Range("H91").Value = 1003.6667
$ws.Range("I91").Value = 1003.6667
$ws.Range("K91").Value = 1003.6667
$ws.Range("M91").Value = 400.3333

# ARM row 92: H92: 98666.336 -> 97666.664; J92: 98666.336 -> 97666.664; L92: 98666.336 -> 97666.664; N92: -103658.336 -> -102658.664
$ws.Range("H92").Value = 97666.664
$ws.Range("J92").Value = 97666.664
$ws.Range("L92").Value = 97666.664
$ws.Range("N92").Value = -102658.664

# ARM row 94: H94: 0 -> 40330; J94: 0 -> 40330; L94: 0 -> 40330; N94: None -> -42132 | add: ['N94']
$ws.Range("H94").Value = 40330
$ws.Range("J94").Value = 40330
$ws.Range("L94").Value = 40330
$ws.Range("N94").Value = -42132

# ARM row 99: H99: 6000 -> 5999; I99: 6000 -> 5999; K99: 6000 -> 5999; M99: -3005 -> -3004
$ws.Range("H99").Value = 5999
$ws.Range("I99").Value = 5999
$ws.Range("K99").Value = 5999
$ws.Range("M99").Value = -3004

$ws = $wb.Worksheets.Item("BSM")
# BSM row 64: H64: 5497.6 -> 4631.3335; I64: 4597 -> 4100; J64: 6098 -> 4897; K64: 4597 -> 4100; L64: 6098 -> 4897; M64: -4372 -> -3875; N64: -6548 -> -5347
$ws.Range("H64").Value = 4631.3335
$ws.Range("I64").Value = 4100
$ws.Range("J64").Value = 4897
$ws.Range("K64").Value = 4100
$ws.Range("L64").Value = 4897
$ws.Range("M64").Value = -3875
$ws.Range("N64").Value = -5347

# BSM row 67: H67: 5497.6 -> 4631.3335; I67: 4597 -> 4100; J67: 6098 -> 4897; K67: 4597 -> 4100; L67: 6098 -> 4897; M67: -3817 -> -3320; N67: -7658 -> -6457
$ws.Range("H67").Value = 4631.3335
$ws.Range("I67").Value = 4100
$ws.Range("J67").Value = 4897
$ws.Range("K67").Value = 4100
$ws.Range("L67").Value = 4897
$ws.Range("M67").Value = -3320
$ws.Range("N67").Value = -6457

# BSM row 86: H86: 1503.5 -> 1928.5; I86: 1000 -> 1850; K86: 1000 -> 1850; M86: 123 -> -727
$ws.Range("H86").Value = 1928.5
$ws.Range("I86").Value = 1850
$ws.Range("K86").Value = 1850
$ws.Range("M86").Value = -727

# BSM row 89: H89: 1503.5 -> 1928.5; I89: 1000 -> 1850; K89: 5000 -> 9250; M89: 616 -> -3634
$ws.Range("H89").Value = 1928.5
$ws.Range("I89").Value = 1850
$ws.Range("K89").Value = 9250
$ws.Range("M89").Value = -3634

# BSM row 92: H92: 39497 -> 38999.5; J92: 39497 -> 38999.5; L92: 39497 -> 38999.5; N92: -44489 -> -43991.5
$ws.Range("H92").Value = 38999.5
$ws.Range("J92").Value = 38999.5
$ws.Range("L92").Value = 38999.5
$ws.Range("N92").Value = -43991.5

# BSM row 96: H96: 48678.25 -> 48489.43; J96: 50571.145 -> 50666.332; L96: 50571.145 -> 50666.332; N96: -56063.145 -> -56158.332
$ws.Range("H96").Value = 48489.43
$ws.Range("J96").Value = 50666.332
$ws.Range("L96").Value = 50666.332
$ws.Range("N96").Value = -56158.332

# BSM row 97: H97: 50000.832 -> 46667.332; I97: 0 -> 29999; J97: 50000.832 -> 50001; K97: 0 -> 29999; L97: 50000.832 -> 50001; N97: -51982.832 -> -51983; M97: None -> -29008 | add: ['M97']
$ws.Range("H97").Value = 46667.332
$ws.Range("I97").Value = 29999
$ws.Range("J97").Value = 50001
$ws.Range("K97").Value = 29999
$ws.Range("L97").Value = 50001
$ws.Range("N97").Value = -51983
$ws.Range("M97").Value = -29008

# BSM row 100: H100: 7939.8 -> 8750; J100: 7939.8 -> 8750; L100: 7939.8 -> 8750; N100: -10103.8 -> -10914
$ws.Range("H100").Value = 8750
$ws.Range("J100").Value = 8750
$ws.Range("L100").Value = 8750
$ws.Range("N100").Value = -10914

$ws = $wb.Worksheets.Item("CRP")
# CRP row 41: H41: 13666 -> 0; I41: 1000 -> 0; J41: 19999 -> 0; K41: 1000 -> 0; L41: 19999 -> 0 | remove: ['M41', 'N41']
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()

# CRP row 62: H62: 4948 -> 4997; I62: 4948 -> 4997; K62: 4948 -> 4997; M62: -4324 -> -4373
$ws.Range("H62").Value = 4997
$ws.Range("I62").Value = 4997
$ws.Range("K62").Value = 4997
$ws.Range("M62").Value = -4373

# CRP row 65: H65: 4948 -> 4997; I65: 4948 -> 4997; K65: 24740 -> 24985; M65: -21620 -> -21865
$ws.Range("H65").Value = 4997
$ws.Range("I65").Value = 4997
$ws.Range("K65").Value = 24985
$ws.Range("M65").Value = -21865

# CRP row 68: H68: 49999 -> 49800; J68: 49999 -> 49800; L68: 49999 -> 49800; N68: -51497 -> -51298
$ws.Range("H68").Value = 49800
$ws.Range("J68").Value = 49800
$ws.Range("L68").Value = 49800
$ws.Range("N68").Value = -51298

# CRP row 71: H71: 49999 -> 49800; J71: 49999 -> 49800; L71: 149997 -> 149400; N71: -157485 -> -156888
$ws.Range("H71").Value = 49800
$ws.Range("J71").Value = 49800
$ws.Range("L71").Value = 149400
$ws.Range("N71").Value = -156888

# CRP row 96: H96: 16874.334 -> 15655.75; J96: 16874.334 -> 15655.75; L96: 16874.334 -> 15655.75; N96: -22366.334 -> -21147.75
$ws.Range("H96").Value = 15655.75
$ws.Range("J96").Value = 15655.75
$ws.Range("L96").Value = 15655.75
$ws.Range("N96").Value = -21147.75

# CRP row 107: H107: 591.5 -> 591.0833; I107: 620.2 -> 619.7; K107: 620.2 -> 619.7; M107: 1299.8 -> 1300.3
$ws.Range("H107").Value = 591.0833
$ws.Range("I107").Value = 619.7
$ws.Range("K107").Value = 619.7
$ws.Range("M107").Value = 1300.3

# CRP row 122: H122: 1147.25 -> 888; I122: 1147.25 -> 0; J122: 0 -> 888; K122: 3441.75 -> 0; L122: 0 -> 2664; N122: None -> -7564 | add: ['N122'] | remove: ['M122']
$ws.Range("H122").Value = 888
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 888
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 2664
$ws.Range("N122").Value = -7564
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# CUL row 4: H4: 58830384 -> 83342200; I4: 928.5714 -> 350; J4: 100011000 -> 125013130; K4: 2785.7142 -> 1050; L4: 300033000 -> 375039390; M4: -2673.7142 -> -938; N4: -300033224 -> -375039614
$ws.Range("H4").Value = 83342200
$ws.Range("I4").Value = 350
$ws.Range("J4").Value = 125013130
$ws.Range("K4").Value = 1050
$ws.Range("L4").Value = 375039390
$ws.Range("M4").Value = -938
$ws.Range("N4").Value = -375039614

# CUL row 7: H7: 19.375 -> 31.2; I7: 32.5 -> 35.5; J7: 6.25 -> 14; K7: 97.5 -> 106.5; L7: 18.75 -> 42; M7: 14.5 -> 5.5; N7: -242.75 -> -266
$ws.Range("H7").Value = 31.2
$ws.Range("I7").Value = 35.5
$ws.Range("J7").Value = 14
$ws.Range("K7").Value = 106.5
$ws.Range("L7").Value = 42
$ws.Range("M7").Value = 5.5
$ws.Range("N7").Value = -266

# CUL row 108: H108: 1763.5 -> 527; I108: 1763.5 -> 527; K108: 5290.5 -> 1581; M108: -2410.5 -> 1299
$ws.Range("H108").Value = 527
$ws.Range("I108").Value = 527
$ws.Range("K108").Value = 1581
$ws.Range("M108").Value = 1299

$ws = $wb.Worksheets.Item("GSM")
# GSM row 80: H80: 4063 -> 3306.5; I80: 3785.8 -> 3150; J80: 4525 -> 3776; K80: 3785.8 -> 3150; L80: 4525 -> 3776; M80: -2787.8 -> -2152; N80: -6521 -> -5772
$ws.Range("H80").Value = 3306.5
$ws.Range("I80").Value = 3150
$ws.Range("J80").Value = 3776
$ws.Range("K80").Value = 3150
$ws.Range("L80").Value = 3776
$ws.Range("M80").Value = -2152
$ws.Range("N80").Value = -5772

# GSM row 83: H83: 4063 -> 3306.5; I83: 3785.8 -> 3150; J83: 4525 -> 3776; K83: 18929 -> 15750; L83: 22625 -> 18880; M83: -13937 -> -10758; N83: -32609 -> -28864
$ws.Range("H83").Value = 3306.5
$ws.Range("I83").Value = 3150
$ws.Range("J83").Value = 3776
$ws.Range("K83").Value = 15750
$ws.Range("L83").Value = 18880
$ws.Range("M83").Value = -10758
$ws.Range("N83").Value = -28864

# GSM row 122: H122: 8250.5 -> 7296.5; I122: 7503.5 -> 6002.3335; J122: 8997.5 -> 8590.666999999999; K122: 22510.5 -> 18007.0005; L122: 26992.5 -> 25772.001; M122: -20060.5 -> -15557.0005; N122: -31892.5 -> -30672.001
$ws.Range("H122").Value = 7296.5
$ws.Range("I122").Value = 6002.3335
$ws.Range("J122").Value = 8590.666999999999
$ws.Range("K122").Value = 18007.0005
$ws.Range("L122").Value = 25772.001
$ws.Range("M122").Value = -15557.0005
$ws.Range("N122").Value = -30672.001

$ws = $wb.Worksheets.Item("LTW")
# LTW row 14: H14: 5001 -> 0; J14: 5001 -> 0; L14: 5001 -> 0 | remove: ['N14']
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# WVR row 14: H14: 1999.5 -> 2000; J14: 1999.5 -> 2000; L14: 1999.5 -> 2000; N14: -2335.5 -> -2336
$ws.Range("H14").Value = 2000
$ws.Range("J14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("N14").Value = -2336

# WVR row 62: H62: 56250.5 -> 10000.333; I62: 8334 -> 10000.333; J62: 200000 -> 0; K62: 8334 -> 10000.333; L62: 200000 -> 0; M62: -7710 -> -9376.333000000001 | remove: ['N62']
$ws.Range("H62").Value = 10000.333
$ws.Range("I62").Value = 10000.333
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 10000.333
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -9376.333000000001
$ws.Range("N62").ClearContents()

# WVR row 65: H65: 56250.5 -> 10000.333; I65: 8334 -> 10000.333; J65: 200000 -> 0; K65: 41670 -> 50001.665; L65: 1000000 -> 0; M65: -38550 -> -46881.665 | remove: ['N65']
$ws.Range("H65").Value = 10000.333
$ws.Range("I65").Value = 10000.333
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 50001.665
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -46881.665
$ws.Range("N65").ClearContents()
